$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(11, 1).Value = 44310
$ws.Cells.Item(11, 2).Value = "1.1.3"
$ws.Cells.Item(11, 3).Value = "Improvements:`n- after successful voting of DFIP #8 the LM-countdown is no longer needed and removed"

$ws.Range("A10:C10").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)

$ws.Rows.Item(11).RowHeight = $ws.Rows.Item(10).RowHeight

$ws.Range("C11").Select()
